$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2083.5
$ws.Range("J17").Value = 2083.5
$ws.Range("L17").Value = 6250.5
$ws.Range("N17").Value = -6586.5
# Row 33
$ws.Range("H33").Value = 119.789474
$ws.Range("I33").Value = 98.77778000000001
$ws.Range("J33").Value = 498
$ws.Range("K33").Value = 98.77778000000001
$ws.Range("L33").Value = 498
$ws.Range("M33").Value = 130.22222
$ws.Range("N33").Value = -956
# Row 40
$ws.Range("H40").Value = 2190
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2190
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2190
$ws.Range("N40").Value = -2540
$ws.Range("M40").ClearContents()
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
# Row 137
$ws.Range("H137").Value = 3576.3333
$ws.Range("I137").Value = 916.2
$ws.Range("K137").Value = 2748.6
$ws.Range("M137").Value = -198.6000000000004

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6434.45
$ws.Range("I32").Value = 6052.4707
$ws.Range("J32").Value = 8599
$ws.Range("K32").Value = 6052.4707
$ws.Range("L32").Value = 8599
$ws.Range("M32").Value = -5765.4707
$ws.Range("N32").Value = -9173
# Row 61
$ws.Range("H61").Value = 1543.3636
$ws.Range("I61").Value = 663.3333
$ws.Range("K61").Value = 663.3333
$ws.Range("M61").Value = -451.3333
# Row 63
$ws.Range("H63").Value = 3723.75
$ws.Range("I63").Value = 3435
$ws.Range("K63").Value = 3435
$ws.Range("M63").Value = -2749
# Row 66
$ws.Range("H66").Value = 3723.75
$ws.Range("I66").Value = 3435
$ws.Range("K66").Value = 17175
$ws.Range("M66").Value = -13743
# Row 132
$ws.Range("H132").Value = 2445
$ws.Range("I132").Value = 1561.9524
$ws.Range("J132").Value = 3535.8235
$ws.Range("K132").Value = 4685.857199999999
$ws.Range("L132").Value = 10607.4705
$ws.Range("M132").Value = -2155.857199999999
$ws.Range("N132").Value = -15667.4705
# Row 136
$ws.Range("H136").Value = 1543.3636
$ws.Range("I136").Value = 663.3333
$ws.Range("K136").Value = 1989.9999
$ws.Range("M136").Value = 560.0001
# Row 139
$ws.Range("H139").Value = 62499.25
$ws.Range("J139").Value = 62499.25
$ws.Range("L139").Value = 62499.25
$ws.Range("N139").Value = -72779.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5181.8
$ws.Range("I86").Value = 4146.2
$ws.Range("J86").Value = 6217.4
$ws.Range("K86").Value = 4146.2
$ws.Range("L86").Value = 6217.4
$ws.Range("M86").Value = -3023.2
$ws.Range("N86").Value = -8463.4
# Row 89
$ws.Range("H89").Value = 5181.8
$ws.Range("I89").Value = 4146.2
$ws.Range("J89").Value = 6217.4
$ws.Range("K89").Value = 20731
$ws.Range("L89").Value = 31087
$ws.Range("M89").Value = -15115
$ws.Range("N89").Value = -42319
# Row 94
$ws.Range("H94").Value = 487.55554
$ws.Range("I94").Value = 450.57144
$ws.Range("K94").Value = 450.57144
$ws.Range("M94").Value = 0.4285600000000045
# Row 105
$ws.Range("H105").Value = 3176.8
$ws.Range("I105").Value = 3085.2222
$ws.Range("K105").Value = 3085.2222
$ws.Range("M105").Value = -1338.2222
# Row 138
$ws.Range("H138").Value = 123324.25
$ws.Range("J138").Value = 123324.25
$ws.Range("L138").Value = 123324.25
$ws.Range("N138").Value = -133604.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3337.8333
$ws.Range("J58").Value = 2756.75
$ws.Range("L58").Value = 2756.75
$ws.Range("N58").Value = -3162.75
# Row 94
$ws.Range("H94").Value = 866.3333
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 299.5
$ws.Range("K94").Value = 2000
$ws.Range("L94").Value = 299.5
$ws.Range("M94").Value = -1549
$ws.Range("N94").Value = -1201.5
# Row 136
$ws.Range("H136").Value = 3337.8333
$ws.Range("J136").Value = 2756.75
$ws.Range("L136").Value = 8270.25
$ws.Range("N136").Value = -13370.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 96.666664
$ws.Range("I33").Value = 17
$ws.Range("J33").Value = 495
$ws.Range("K33").Value = 102
$ws.Range("L33").Value = 2970
$ws.Range("M33").Value = 181
$ws.Range("N33").Value = -3536
# Row 60
$ws.Range("H60").Value = 1000
$ws.Range("I60").Value = 1000
$ws.Range("K60").Value = 3000
$ws.Range("M60").Value = -2749
# Row 86
$ws.Range("H86").Value = 312.5
$ws.Range("J86").Value = 250
$ws.Range("L86").Value = 750
$ws.Range("N86").Value = -3122
# Row 89
$ws.Range("H89").Value = 312.5
$ws.Range("J89").Value = 250
$ws.Range("L89").Value = 2250
$ws.Range("N89").Value = -14106

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3632.2856
$ws.Range("I102").Value = 2128
$ws.Range("J102").Value = 5638
$ws.Range("K102").Value = 2128
$ws.Range("L102").Value = 5638
$ws.Range("M102").Value = -506
$ws.Range("N102").Value = -8882
# Row 132
$ws.Range("H132").Value = 3399.625
$ws.Range("I132").Value = 2100
$ws.Range("J132").Value = 4179.4
$ws.Range("K132").Value = 6300
$ws.Range("L132").Value = 12538.2
$ws.Range("M132").Value = -3770
$ws.Range("N132").Value = -17598.2

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 400
$ws.Range("K16").Value = 400
$ws.Range("M16").Value = -230
# Row 82
$ws.Range("H82").Value = 1300
$ws.Range("I82").Value = 1200
$ws.Range("J82").Value = 1400
$ws.Range("K82").Value = 1200
$ws.Range("L82").Value = 1400
$ws.Range("M82").Value = -839
$ws.Range("N82").Value = -2122
# Row 85
$ws.Range("H85").Value = 1300
$ws.Range("I85").Value = 1200
$ws.Range("J85").Value = 1400
$ws.Range("K85").Value = 1200
$ws.Range("L85").Value = 1400
$ws.Range("M85").Value = 48
$ws.Range("N85").Value = -3896
# Row 122
$ws.Range("H122").Value = 6749.5
$ws.Range("I122").Value = 6999
$ws.Range("J122").Value = 6666.3335
$ws.Range("K122").Value = 20997
$ws.Range("L122").Value = 19999.0005
$ws.Range("M122").Value = -18547
$ws.Range("N122").Value = -24899.0005

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 848.5
$ws.Range("I107").Value = 848.5
$ws.Range("K107").Value = 2545.5
$ws.Range("M107").Value = -625.5
# Row 113
$ws.Range("H113").Value = 852.625
$ws.Range("I113").Value = 852.625
$ws.Range("K113").Value = 2557.875
$ws.Range("M113").Value = -387.875
# Row 132
$ws.Range("H132").Value = 3522
$ws.Range("I132").Value = 2901.5557
$ws.Range("K132").Value = 8704.667099999999
$ws.Range("M132").Value = -6174.667099999999
